# Adds header labels to the "Cenarios_PNE" pipeline constant sheets and
# fixes the accented Portuguese labels (typo fixes), as described in commit
# "Colocando header nos gráficos".

$wb = $excel.ActiveWorkbook

# xlPasteFormats constant used with PasteSpecial below.
$xlPasteFormats = -4122

# --- Sheets 1-4 share an identical layout: add "Fonte/Tecnologia" header in A1
#     (copying the header style already used by B1:E1) and correct the
#     accented labels in column A (rows 2-12), clearing their bold/border
#     style in the process. ---
$sheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

$labelFixes = @{
    "Gas Natural" = "Gás Natural"
    "Carvao"      = "Carvão"
    "Oleos Comb"  = "Óleos Comb"
    "Eolica"      = "Eólica"
    "Pot Compl"   = "Pot. Compl."
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Add header cell, copying the style used by the rest of row 1.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial($xlPasteFormats)
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    for ($r = 2; $r -le 12; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value()
        if ($labelFixes.ContainsKey($current)) {
            $cell.Value = $labelFixes[$current]
        }
        $cell.Style = "Normal"
    }
}

# --- Sheet 5: "Emissoes Totais (MtCO2eq)" ---
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial($xlPasteFormats)
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A2").Style = "Normal"

$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A3").Style = "Normal"

# Remove the "Teto" row entirely.
$ws5.Rows.Item(4).Delete()

# --- Sheet 6: "Custo Total (bilhões de R$)" ---
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial($xlPasteFormats)
$ws6.Range("A1").Value = "Tipo Expansão"

# B1's label changes from "Custo" to the (text) year "2015"; reuse the
# existing "2015" text cell from another sheet so the value stays a string
# (t="inlineStr") instead of becoming a number, while keeping B1's style.
$wb.Worksheets.Item("Potencia Acumulada - SIN (MW)").Range("B1").Copy()
$ws6.Range("B1").PasteSpecial(-4104)  # xlPasteAll

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("A2").Style = "Normal"
$ws6.Range("B2").Value = 608

$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("A3").Style = "Normal"
$ws6.Range("B3").Value = 99
